$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the mastery criterion value in J2 from 4 to 8
$ws.Range("J2").Value = 8

# Update the active cell selection to J3
$ws.Range("J3").Select()
